$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(16, 8).Value = 29999  # H16
$ws.Cells.Item(16, 10).Value = 29999  # J16
$ws.Cells.Item(16, 12).Value = 29999  # L16
$ws.Cells.Item(16, 14).Value = -30459  # N16

$ws.Cells.Item(33, 8).Value = 134.85715  # H33
$ws.Cells.Item(33, 9).Value = 139.44444  # I33
$ws.Cells.Item(33, 10).Value = 107.333336  # J33
$ws.Cells.Item(33, 11).Value = 139.44444  # K33
$ws.Cells.Item(33, 12).Value = 107.333336  # L33
$ws.Cells.Item(33, 13).Value = 89.55556000000001  # M33
$ws.Cells.Item(33, 14).Value = -565.333336  # N33

$ws.Cells.Item(40, 8).Value = 5166.3335  # H40
$ws.Cells.Item(40, 9).Value = 4500  # I40
$ws.Cells.Item(40, 10).Value = 5499.5  # J40
$ws.Cells.Item(40, 11).Value = 4500  # K40
$ws.Cells.Item(40, 12).Value = 5499.5  # L40
$ws.Cells.Item(40, 13).Value = -4325  # M40
$ws.Cells.Item(40, 14).Value = -5849.5  # N40

$ws.Cells.Item(51, 8).Value = 12558.833  # H51
$ws.Cells.Item(51, 10).Value = 7070.6  # J51
$ws.Cells.Item(51, 12).Value = 7070.6  # L51
$ws.Cells.Item(51, 14).Value = -8038.6  # N51

$ws.Cells.Item(74, 8).Value = 5929.24  # H74
$ws.Cells.Item(74, 9).Value = 5388.7856  # I74
$ws.Cells.Item(74, 11).Value = 5388.7856  # K74
$ws.Cells.Item(74, 13).Value = -4452.7856  # M74

$ws.Cells.Item(75, 8).Value = 54999.5  # H75
$ws.Cells.Item(75, 10).Value = 54999.5  # J75
$ws.Cells.Item(75, 12).Value = 54999.5  # L75
$ws.Cells.Item(75, 14).Value = -56871.5  # N75

$ws.Cells.Item(76, 8).Value = 4453.5713  # H76
$ws.Cells.Item(76, 9).Value = 4200  # I76
$ws.Cells.Item(76, 10).Value = 4555  # J76
$ws.Cells.Item(76, 11).Value = 4200  # K76
$ws.Cells.Item(76, 12).Value = 4555  # L76
$ws.Cells.Item(76, 13).Value = -3885  # M76
$ws.Cells.Item(76, 14).Value = -5185  # N76

$ws.Cells.Item(77, 8).Value = 5929.24  # H77
$ws.Cells.Item(77, 9).Value = 5388.7856  # I77
$ws.Cells.Item(77, 11).Value = 26943.928  # K77
$ws.Cells.Item(77, 13).Value = -22263.928  # M77

$ws.Cells.Item(78, 8).Value = 54999.5  # H78
$ws.Cells.Item(78, 10).Value = 54999.5  # J78
$ws.Cells.Item(78, 12).Value = 164998.5  # L78
$ws.Cells.Item(78, 14).Value = -174358.5  # N78

$ws.Cells.Item(79, 8).Value = 4453.5713  # H79
$ws.Cells.Item(79, 9).Value = 4200  # I79
$ws.Cells.Item(79, 10).Value = 4555  # J79
$ws.Cells.Item(79, 11).Value = 4200  # K79
$ws.Cells.Item(79, 12).Value = 4555  # L79
$ws.Cells.Item(79, 13).Value = -3108  # M79
$ws.Cells.Item(79, 14).Value = -6739  # N79

$ws.Cells.Item(112, 8).Value = 3232.7354  # H112
$ws.Cells.Item(112, 9).Value = 869.2  # I112
$ws.Cells.Item(112, 10).Value = 4217.5415  # J112
$ws.Cells.Item(112, 11).Value = 2607.6  # K112
$ws.Cells.Item(112, 12).Value = 12652.6245  # L112
$ws.Cells.Item(112, 13).Value = -1499.6  # M112
$ws.Cells.Item(112, 14).Value = -14868.6245  # N112

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(132, 8).Value = 1768.0571  # H132
$ws.Cells.Item(132, 9).Value = 1571.1724  # I132
$ws.Cells.Item(132, 10).Value = 2719.6667  # J132
$ws.Cells.Item(132, 11).Value = 4713.5172  # K132
$ws.Cells.Item(132, 12).Value = 8159.000100000001  # L132
$ws.Cells.Item(132, 13).Value = -2183.5172  # M132
$ws.Cells.Item(132, 14).Value = -13219.0001  # N132

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(70, 8).Value = 167222  # H70
$ws.Cells.Item(70, 10).Value = 167222  # J70
$ws.Cells.Item(70, 12).Value = 167222  # L70
$ws.Cells.Item(70, 14).Value = -167808  # N70

$ws.Cells.Item(73, 8).Value = 167222  # H73
$ws.Cells.Item(73, 10).Value = 167222  # J73
$ws.Cells.Item(73, 12).Value = 167222  # L73
$ws.Cells.Item(73, 14).Value = -169250  # N73

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(4, 8).Value = 0  # H4
$ws.Cells.Item(4, 10).Value = 0  # J4
$ws.Cells.Item(4, 12).Value = 0  # L4
$ws.Cells.Item(4, 14).ClearContents()  # N4

$ws.Cells.Item(7, 8).Value = 132.5  # H7
$ws.Cells.Item(7, 9).Value = 145.57143  # I7
$ws.Cells.Item(7, 11).Value = 145.57143  # K7
$ws.Cells.Item(7, 13).Value = -32.57142999999999  # M7

$ws.Cells.Item(48, 8).Value = 50000  # H48
$ws.Cells.Item(48, 10).Value = 50000  # J48
$ws.Cells.Item(48, 12).Value = 50000  # L48
$ws.Cells.Item(48, 14).Value = -50952  # N48

$ws.Cells.Item(132, 8).Value = 3833.8125  # H132
$ws.Cells.Item(132, 9).Value = 1781.7916  # I132
$ws.Cells.Item(132, 10).Value = 9989.875  # J132
$ws.Cells.Item(132, 11).Value = 5345.3748  # K132
$ws.Cells.Item(132, 12).Value = 29969.625  # L132
$ws.Cells.Item(132, 13).Value = -2815.3748  # M132
$ws.Cells.Item(132, 14).Value = -35029.625  # N132

$ws.Cells.Item(140, 8).Value = 114861.5  # H140
$ws.Cells.Item(140, 10).Value = 114861.5  # J140
$ws.Cells.Item(140, 12).Value = 114861.5  # L140
$ws.Cells.Item(140, 14).Value = -125221.5  # N140

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(34, 8).Value = 4074.3333  # H34
$ws.Cells.Item(34, 9).Value = 112.166664  # I34
$ws.Cells.Item(34, 10).Value = 11998.667  # J34
$ws.Cells.Item(34, 11).Value = 336.499992  # K34
$ws.Cells.Item(34, 12).Value = 35996.001  # L34
$ws.Cells.Item(34, 13).Value = -252.499992  # M34
$ws.Cells.Item(34, 14).Value = -36164.001  # N34

$ws.Cells.Item(39, 8).Value = 6355.1113  # H39
$ws.Cells.Item(39, 10).Value = 8113.7144  # J39
$ws.Cells.Item(39, 12).Value = 24341.1432  # L39
$ws.Cells.Item(39, 14).Value = -24929.1432  # N39

$ws.Cells.Item(55, 8).Value = 1193  # H55
$ws.Cells.Item(55, 10).Value = 2000  # J55
$ws.Cells.Item(55, 12).Value = 6000  # L55
$ws.Cells.Item(55, 14).Value = -6354  # N55

$ws.Cells.Item(114, 8).Value = 8000  # H114
$ws.Cells.Item(114, 9).Value = 0  # I114
$ws.Cells.Item(114, 11).Value = 0  # K114
$ws.Cells.Item(114, 13).ClearContents()  # M114

$ws.Cells.Item(130, 8).Value = 23329.9  # H130
$ws.Cells.Item(130, 10).Value = 24977.666  # J130
$ws.Cells.Item(130, 12).Value = 74932.99800000001  # L130
$ws.Cells.Item(130, 14).Value = -84972.99800000001  # N130

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(33, 8).Value = 5000  # H33
$ws.Cells.Item(33, 10).Value = 5000  # J33
$ws.Cells.Item(33, 12).Value = 5000  # L33
$ws.Cells.Item(33, 14).Value = -5504  # N33

$ws.Cells.Item(47, 8).Value = 36779.5  # H47
$ws.Cells.Item(47, 9).Value = 0  # I47
$ws.Cells.Item(47, 10).Value = 36779.5  # J47
$ws.Cells.Item(47, 11).Value = 0  # K47
$ws.Cells.Item(47, 12).Value = 36779.5  # L47
$ws.Cells.Item(47, 13).ClearContents()  # M47
$ws.Cells.Item(47, 14).Value = -37915.5  # N47

$ws.Cells.Item(48, 8).Value = 0  # H48
$ws.Cells.Item(48, 10).Value = 0  # J48
$ws.Cells.Item(48, 12).Value = 0  # L48
$ws.Cells.Item(48, 14).ClearContents()  # N48

$ws.Cells.Item(70, 8).Value = 9891.143  # H70
$ws.Cells.Item(70, 9).Value = 11309.75  # I70
$ws.Cells.Item(70, 10).Value = 7999.6665  # J70
$ws.Cells.Item(70, 11).Value = 11309.75  # K70
$ws.Cells.Item(70, 12).Value = 7999.6665  # L70
$ws.Cells.Item(70, 13).Value = -11039.75  # M70
$ws.Cells.Item(70, 14).Value = -8539.666499999999  # N70

$ws.Cells.Item(73, 8).Value = 9891.143  # H73
$ws.Cells.Item(73, 9).Value = 11309.75  # I73
$ws.Cells.Item(73, 10).Value = 7999.6665  # J73
$ws.Cells.Item(73, 11).Value = 11309.75  # K73
$ws.Cells.Item(73, 12).Value = 7999.6665  # L73
$ws.Cells.Item(73, 13).Value = -10373.75  # M73
$ws.Cells.Item(73, 14).Value = -9871.666499999999  # N73

$ws.Cells.Item(80, 8).Value = 2526.125  # H80
$ws.Cells.Item(80, 10).Value = 1968.6666  # J80
$ws.Cells.Item(80, 12).Value = 1968.6666  # L80
$ws.Cells.Item(80, 14).Value = -3964.6666  # N80

$ws.Cells.Item(83, 8).Value = 2526.125  # H83
$ws.Cells.Item(83, 10).Value = 1968.6666  # J83
$ws.Cells.Item(83, 12).Value = 9843.333000000001  # L83
$ws.Cells.Item(83, 14).Value = -19827.333  # N83

$ws.Cells.Item(97, 8).Value = 2936.12  # H97
$ws.Cells.Item(97, 9).Value = 528  # I97
$ws.Cells.Item(97, 10).Value = 5544.9165  # J97
$ws.Cells.Item(97, 11).Value = 528  # K97
$ws.Cells.Item(97, 12).Value = 5544.9165  # L97
$ws.Cells.Item(97, 13).Value = -32  # M97
$ws.Cells.Item(97, 14).Value = -6536.9165  # N97

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(3, 8).Value = 20000  # H3
$ws.Cells.Item(3, 10).Value = 0  # J3
$ws.Cells.Item(3, 12).Value = 0  # L3
$ws.Cells.Item(3, 14).ClearContents()  # N3

$ws.Cells.Item(15, 8).Value = 20000  # H15
$ws.Cells.Item(15, 10).Value = 0  # J15
$ws.Cells.Item(15, 12).Value = 0  # L15
$ws.Cells.Item(15, 14).ClearContents()  # N15

$ws.Cells.Item(40, 8).Value = 2455.9285  # H40
$ws.Cells.Item(40, 9).Value = 2158.7  # I40
$ws.Cells.Item(40, 11).Value = 2158.7  # K40
$ws.Cells.Item(40, 13).Value = -2022.7  # M40

$ws.Cells.Item(82, 8).Value = 1360.6875  # H82
$ws.Cells.Item(82, 10).Value = 1485.3334  # J82
$ws.Cells.Item(82, 12).Value = 1485.3334  # L82
$ws.Cells.Item(82, 14).Value = -2207.3334  # N82

$ws.Cells.Item(85, 8).Value = 1360.6875  # H85
$ws.Cells.Item(85, 10).Value = 1485.3334  # J85
$ws.Cells.Item(85, 12).Value = 1485.3334  # L85
$ws.Cells.Item(85, 14).Value = -3981.3334  # N85

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(19, 8).Value = 0  # H19
$ws.Cells.Item(19, 10).Value = 0  # J19
$ws.Cells.Item(19, 12).Value = 0  # L19
$ws.Cells.Item(19, 14).ClearContents()  # N19

$ws.Cells.Item(81, 8).Value = 4335.364  # H81
$ws.Cells.Item(81, 9).Value = 5799.7334  # I81
$ws.Cells.Item(81, 10).Value = 1197.4286  # J81
$ws.Cells.Item(81, 11).Value = 11599.4668  # K81
$ws.Cells.Item(81, 12).Value = 2394.8572  # L81
$ws.Cells.Item(81, 13).Value = -10538.4668  # M81
$ws.Cells.Item(81, 14).Value = -4516.8572  # N81

$ws.Cells.Item(84, 8).Value = 4335.364  # H84
$ws.Cells.Item(84, 9).Value = 5799.7334  # I84
$ws.Cells.Item(84, 10).Value = 1197.4286  # J84
$ws.Cells.Item(84, 11).Value = 57997.334  # K84
$ws.Cells.Item(84, 12).Value = 11974.286  # L84
$ws.Cells.Item(84, 13).Value = -52693.334  # M84
$ws.Cells.Item(84, 14).Value = -22582.286  # N84

$ws.Cells.Item(103, 8).Value = 35520.4  # H103
$ws.Cells.Item(103, 10).Value = 35520.4  # J103
$ws.Cells.Item(103, 12).Value = 35520.4  # L103
$ws.Cells.Item(103, 14).Value = -37864.4  # N103

$ws.Cells.Item(137, 8).Value = 62023.4  # H137
$ws.Cells.Item(137, 10).Value = 62023.4  # J137
$ws.Cells.Item(137, 12).Value = 62023.4  # L137
$ws.Cells.Item(137, 14).Value = -72223.39999999999  # N137

$ws.Cells.Item(141, 8).Value = 93496.586  # H141
$ws.Cells.Item(141, 10).Value = 93496.586  # J141
$ws.Cells.Item(141, 12).Value = 93496.586  # L141
$ws.Cells.Item(141, 14).Value = -103856.586  # N141
